# Generate Report for Handoff
#
# The localization status report was regenerated: a new handoff id
# (f3ae1452-70bb-4ff0-848b-a5c7e415c328) replaces the previous one
# (acbe35b7-8549-4ca7-9547-92bfdf29c244), new xlf content-hashes are
# baked into the target filenames, and the handoff timestamps advance.
#
# This updates, on every worksheet, the cell text that shows those
# filenames/timestamps and keeps each cell's existing hyperlink (same
# target URL / relationship) pointing at it, just with the refreshed
# display text - mirroring what the CI job that produces this report
# would do on a fresh run.

$wb = $excel.ActiveWorkbook

$oldId = "acbe35b7-8549-4ca7-9547-92bfdf29c244"
$newId = "f3ae1452-70bb-4ff0-848b-a5c7e415c328"

$oldZhHash = "0c1534f788e0cc17161d6d1ab5f741abe753040f"
$newZhHash = "7f3e821fce2dfd45848989f58fa163f2d10736d7"

$oldDeHash = "0c1534f788e0cc17161d6d1ab5f741abe753040f"
$newDeHash = "7f3e821fce2dfd45848989f58fa163f2d10736d7"

$newMdName = "$newId.md"
$newZhXlf  = "$newId.$newZhHash.zh-cn.xlf"
$newDeXlf  = "$newId.$newDeHash.de-de.xlf"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = "2016-03-24 13:18:49"

$mdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/a85fa65b5bed5db82ad658e3f0f15bb5d783dcaa/e2e/$oldId.md"

# The engine's Hyperlinks collection only supports clearing *all*
# hyperlinks on a sheet at once, so rebuild every hyperlink on the
# sheet (same target addresses, refreshed display text).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdTarget, "", "", $newMdName)

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhXlf
$wsZh.Range("E2").Value = "2016-03-24 13:18:45"

$zhXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e966494967337814b615f369cfba4d01cdd396f9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldId.$oldZhHash.zh-cn.xlf"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdTarget, "", "", $newMdName)
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlfTarget, "", "", $newZhXlf)

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeXlf
$wsDe.Range("E2").Value = "2016-03-24 13:18:49"

$deXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ec225aeb014af991381aa1fa7f45ccf80ae07cd4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldId.$oldDeHash.de-de.xlf"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdTarget, "", "", $newMdName)
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlfTarget, "", "", $newDeXlf)
